$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ser" row (row 3) — rows below shift up
$ws.Rows.Item(3).Delete()

# Add a new "p1" row at the bottom (row 5) with its scores
$ws.Range("A5").Value = "p1"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 35
